$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 updates -------------------------------------------------------

# D4: same "Actual" date label already used in D3 ("2025.01.12").
# Assigning the literal string directly would make Excel's auto-detection
# turn it into a date serial number, so copy the existing text cell
# (C3, which already holds that exact shared string) and paste its value
# into D4. This keeps D4 as plain text using the existing shared string,
# matching how D3 was authored, and avoids creating any new cell style.
$ws.Range("C3").Copy() | Out-Null
$ws.Range("D4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$excel.CutCopyMode = $false

# E4: progress goes from 0% to 90% (keeps its existing percentage style).
$ws.Range("E4").Value = 0.9

# F4: new Subject value for this session.
$ws.Range("F4").Value = "K-NN & Parzen Window"

# --- Column F sizing -------------------------------------------------------
# New column needed for the Subject text; approximate the bestFit width.
$ws.Columns("F").ColumnWidth = 19

# --- Selection change --------------------------------------------------
$ws.Range("F4").Select()
